$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.696344631456725
$ws.Range("C2").Value = 0.2111826404313319
$ws.Range("E2").Value = 0.1202738919542483
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.4402459494588271
$ws.Range("H2").Value = 0.5399556131438601
$ws.Range("I2").Value = 0.3815909797095571
$ws.Range("M2").Value = 0.3286786654464038
$ws.Range("N2").Value = 0.9353039029448524

$ws.Range("B3").Value = 0.6098563210493353
$ws.Range("C3").Value = 0.1846755514988843
$ws.Range("E3").Value = 0.1143065586289254
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.4265770218188862
$ws.Range("H3").Value = 0.5388796194290819
$ws.Range("I3").Value = 0.3829905893173908
$ws.Range("M3").Value = 0.2927292755368498
$ws.Range("N3").Value = 0.9490907083500488

$ws.Range("B4").Value = 0.5567321126492004
$ws.Range("C4").Value = 0.1683439874708768
$ws.Range("E4").Value = 0.1107395920263698
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.4186560171176694
$ws.Range("H4").Value = 0.5386230858463108
$ws.Range("I4").Value = 0.3841979447782577
$ws.Range("M4").Value = 0.2707506908111057
$ws.Range("N4").Value = 0.9580302143529984

$ws.Range("B5").Value = 0.5350788195066798
$ws.Range("C5").Value = 0.1616744999048478
$ws.Range("E5").Value = 0.1093101728013579
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.4155458298511121
$ws.Range("H5").Value = 0.5386198408775869
$ws.Range("I5").Value = 0.3847770598594877
$ws.Range("M5").Value = 0.261817687103381
$ws.Range("N5").Value = 0.961792291447086

$ws.Range("B6").Value = 0.531483028821782
$ws.Range("C6").Value = 0.1605661724988465
$ws.Range("E6").Value = 0.1090742702730054
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.4150364660140582
$ws.Range("H6").Value = 0.5386254104534629
$ws.Range("I6").Value = 0.3848784709887134
$ws.Range("M6").Value = 0.2603357758677447
$ws.Range("N6").Value = 0.9624241746706197

$ws.Range("B7").Value = 0.5564401075912997
$ws.Range("C7").Value = 0.1682540982702676
$ws.Range("E7").Value = 0.1107202169216706
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.4186135967871394
$ws.Range("H7").Value = 0.5386226323828822
$ws.Range("I7").Value = 0.38420540275796
$ws.Range("M7").Value = 0.2706301226651533
$ws.Range("N7").Value = 0.9580804687498841

$ws.Range("B8").Value = 0.6665277666840836
$ws.Range("C8").Value = 0.2020545512467038
$ws.Range("E8").Value = 0.1181960744998349
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.4354344126698777
$ws.Range("H8").Value = 0.5395005335414567
$ws.Range("I8").Value = 0.3820010998871339
$ws.Range("M8").Value = 0.3162635058303209
$ws.Range("N8").Value = 0.9399590674097169

$ws.Range("B9").Value = 0.8822486809775683
$ws.Range("C9").Value = 0.2679012645247951
$ws.Range("E9").Value = 0.1336377011825789
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.4722081823344553
$ws.Range("H9").Value = 0.5444444606195731
$ws.Range("I9").Value = 0.3804579671202291
$ws.Range("M9").Value = 0.4065188906039623
$ws.Range("N9").Value = 0.9081920617540788

$ws.Range("B10").Value = 1.040655127644925
$ws.Range("C10").Value = 0.3160306151322345
$ws.Range("E10").Value = 0.1454767503955026
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.5016007213241949
$ws.Range("H10").Value = 0.5500644729294066
$ws.Range("I10").Value = 0.3810447542030744
$ws.Range("M10").Value = 0.4733307572301158
$ws.Range("N10").Value = 0.8871572875003508

$ws.Range("B11").Value = 1.11270432345674
$ws.Range("C11").Value = 0.3378760758614874
$ws.Range("E11").Value = 0.1509738323944063
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.5155021557538362
$ws.Range("H11").Value = 0.553057914211692
$ws.Range("I11").Value = 0.3816909756147879
$ws.Range("M11").Value = 0.5038415621186516
$ws.Range("N11").Value = 0.878089951822318

$ws.Range("B12").Value = 1.139985928716158
$ws.Range("C12").Value = 0.3461415927382347
$ws.Range("E12").Value = 0.1530717519622868
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.5208436667198697
$ws.Range("H12").Value = 0.5542546696412529
$ws.Range("I12").Value = 0.3819906778924889
$ws.Range("M12").Value = 0.5154126086090542
$ws.Range("N12").Value = 0.8747286637057421

$ws.Range("B13").Value = 1.134110434751449
$ws.Range("C13").Value = 0.3443617698109449
$ws.Range("E13").Value = 0.1526191984618208
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.5196898205137046
$ws.Range("H13").Value = 0.5539941097957666
$ws.Range("I13").Value = 0.3819236786457623
$ws.Range("M13").Value = 0.5129198033905027
$ws.Range("N13").Value = 0.8754493577315934

$ws.Range("B14").Value = 1.114948836673364
$ws.Range("C14").Value = 0.3385562221074849
$ws.Range("E14").Value = 0.1511461013960513
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.5159400491241257
$ws.Range("H14").Value = 0.5531551030925357
$ws.Range("I14").Value = 0.3817145275846059
$ws.Range("M14").Value = 0.5047931710909239
$ws.Range("N14").Value = 0.8778119661470249

$ws.Range("B15").Value = 1.103211552675077
$ws.Range("C15").Value = 0.3349992639528239
$ws.Range("E15").Value = 0.1502459168824331
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.5136533085975827
$ws.Range("H15").Value = 0.5526494297217397
$ws.Range("I15").Value = 0.381593591553802
$ws.Range("M15").Value = 0.4998176328174537
$ws.Range("N15").Value = 0.8792685556702615

$ws.Range("B16").Value = 1.035946300211776
$ws.Range("C16").Value = 0.3146019919931859
$ws.Range("E16").Value = 0.1451197685528385
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.5007029958753009
$ws.Range("H16").Value = 0.5498776680991853
$ws.Range("I16").Value = 0.3810101928809999
$ws.Range("M16").Value = 0.471339202304037
$ws.Range("N16").Value = 0.8877599675458363

$ws.Range("B17").Value = 0.9946782748679652
$ws.Range("C17").Value = 0.3020765172999518
$ws.Range("E17").Value = 0.1420037755997328
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.4928950396395209
$ws.Range("H17").Value = 0.5482894464335573
$ws.Range("I17").Value = 0.3807497619340836
$ws.Range("M17").Value = 0.4538990029109726
$ws.Range("N17").Value = 0.8930977585530293

$ws.Range("B18").Value = 0.9709410107109306
$ws.Range("C18").Value = 0.294867591885918
$ws.Range("E18").Value = 0.1402220118284916
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.488453992048278
$ws.Range("H18").Value = 0.5474170404857972
$ws.Range("I18").Value = 0.3806356423251458
$ws.Range("M18").Value = 0.4438789245604795
$ws.Range("N18").Value = 0.8962151047850995

$ws.Range("B19").Value = 0.9629038213453782
$ws.Range("C19").Value = 0.2924259788357233
$ws.Range("E19").Value = 0.1396205278065494
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.4869588654577655
$ws.Range("H19").Value = 0.5471287056817857
$ws.Range("I19").Value = 0.3806031166209394
$ws.Range("M19").Value = 0.4404881881824565
$ws.Range("N19").Value = 0.8972786835494375

$ws.Range("B20").Value = 0.9990714311950342
$ws.Range("C20").Value = 0.3034103509565114
$ws.Range("E20").Value = 0.1423343925160836
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.4937210397792313
$ws.Range("H20").Value = 0.5484542591719475
$ws.Range("I20").Value = 0.3807737899285613
$ws.Range("M20").Value = 0.4557543945223728
$ws.Range("N20").Value = 0.8925246561493445

$ws.Range("B21").Value = 1.120577115573951
$ws.Range("C21").Value = 0.3402616379575534
$ws.Range("E21").Value = 0.1515783413447096
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.5170393406228015
$ws.Range("H21").Value = 0.5533998212905402
$ws.Range("I21").Value = 0.3817744642428096
$ws.Range("M21").Value = 0.5071796886834221
$ws.Range("N21").Value = 0.8771160465361554

$ws.Range("B22").Value = 1.19997710976935
$ws.Range("C22").Value = 0.3643060183462126
$ws.Range("E22").Value = 0.1577148961829522
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.532730475885046
$ws.Range("H22").Value = 0.5570006140681585
$ws.Range("I22").Value = 0.3827492104490346
$ws.Range("M22").Value = 0.5408898943759652
$ws.Range("N22").Value = 0.867467235106723

$ws.Range("B23").Value = 1.157600949375308
$ws.Range("C23").Value = 0.3514767004772352
$ws.Range("E23").Value = 0.1544309102276813
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.5243141860800762
$ws.Range("H23").Value = 0.5550449506811646
$ws.Range("I23").Value = 0.3821994720853041
$ws.Range("M23").Value = 0.5228887878975428
$ws.Range("N23").Value = 0.8725783422292999

$ws.Range("B24").Value = 0.9970853217298554
$ws.Range("C24").Value = 0.3028073492718306
$ws.Range("E24").Value = 0.1421848905349563
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.4933474561433684
$ws.Range("H24").Value = 0.5483796206265481
$ws.Range("I24").Value = 0.3807628159922771
$ws.Range("M24").Value = 0.4549155516857297
$ws.Range("N24").Value = 0.8927836044289634

$ws.Range("B25").Value = 0.8239062183641295
$ws.Range("C25").Value = 0.2501328496340136
$ws.Range("E25").Value = 0.1293748492246962
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.4618475985685393
$ws.Range("H25").Value = 0.5427596258246012
$ws.Range("I25").Value = 0.3805751736322094
$ws.Range("M25").Value = 0.3820167798509928
$ws.Range("N25").Value = 0.9163817050226513
